$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.404023666666667
$ws.Range("H2").Value = 4.212071
$ws.Range("I2").Value = 0.004814659541656092
$ws.Range("J2").Value = 0.004814659541656092
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.33695966666667
$ws.Range("N2").Value = 139.010879
$ws.Range("O2").Value = 0.1993490803952133
$ws.Range("P2").Value = 0.1993490803952133
$ws.Range("Q2").Value = 65.05818801337877
$ws.Range("R2").Value = 585.5236921204089
$ws.Range("S2").Value = 0.0009597979520451813
$ws.Range("T2").Value = 0.0009597979520451811
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.404023666666667
$ws.Range("H3").Value = 4.212071
$ws.Range("I3").Value = 0.004814659541656092
$ws.Range("J3").Value = 0.004814659541656092
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 84.50960033333332
$ws.Range("N3").Value = 253.528801
$ws.Range("O3").Value = 0.3635739425333109
$ws.Range("P3").Value = 0.3635739425333109
$ws.Range("Q3").Value = 118.6534789285412
$ws.Range("R3").Value = 1067.881310356871
$ws.Range("S3").Value = 0.001750484751515529
$ws.Range("T3").Value = 0.001750484751515529
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.404023666666667
$ws.Range("H4").Value = 4.212071
$ws.Range("I4").Value = 0.004814659541656092
$ws.Range("J4").Value = 0.004814659541656092
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 72.52790466666666
$ws.Range("N4").Value = 217.583714
$ws.Range("O4").Value = 0.3120267536390091
$ws.Range("P4").Value = 0.3120267536390091
$ws.Range("Q4").Value = 101.8308946457438
$ws.Range("R4").Value = 916.4780518116939
$ws.Range("S4").Value = 0.00150230258666003
$ws.Range("T4").Value = 0.00150230258666003
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.404023666666667
$ws.Range("H5").Value = 4.212071
$ws.Range("I5").Value = 0.004814659541656092
$ws.Range("J5").Value = 0.004814659541656092
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.06683666666666
$ws.Range("N5").Value = 87.20050999999998
$ws.Range("O5").Value = 0.1250502234324667
$ws.Range("P5").Value = 0.1250502234324667
$ws.Range("Q5").Value = 40.81052659513443
$ws.Range("R5").Value = 367.2947393562099
$ws.Range("S5").Value = 0.0006020742514353522
$ws.Range("T5").Value = 0.0006020742514353522
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 21.406858
$ws.Range("H6").Value = 64.220574
$ws.Range("I6").Value = 0.07340811666748523
$ws.Range("J6").Value = 0.07340811666748523
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 46.33695966666667
$ws.Range("N6").Value = 139.010879
$ws.Range("O6").Value = 0.1993490803952133
$ws.Range("P6").Value = 0.1993490803952133
$ws.Range("Q6").Value = 991.9287157360607
$ws.Range("R6").Value = 8927.358441624545
$ws.Range("S6").Value = 0.01463384055120771
$ws.Range("T6").Value = 0.01463384055120771
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 21.406858
$ws.Range("H7").Value = 64.220574
$ws.Range("I7").Value = 0.07340811666748523
$ws.Range("J7").Value = 0.07340811666748523
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 84.50960033333332
$ws.Range("N7").Value = 253.528801
$ws.Range("O7").Value = 0.3635739425333109
$ws.Range("P7").Value = 0.3635739425333109
$ws.Range("Q7").Value = 1809.085013972419
$ws.Range("R7").Value = 16281.76512575177
$ws.Range("S7").Value = 0.02668927839074286
$ws.Range("T7").Value = 0.02668927839074286
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.406858
$ws.Range("H8").Value = 64.220574
$ws.Range("I8").Value = 0.07340811666748523
$ws.Range("J8").Value = 0.07340811666748523
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 72.52790466666666
$ws.Range("N8").Value = 217.583714
$ws.Range("O8").Value = 0.3120267536390091
$ws.Range("P8").Value = 0.3120267536390091
$ws.Range("Q8").Value = 1552.59455623687
$ws.Range("R8").Value = 13973.35100613184
$ws.Range("S8").Value = 0.02290529633450905
$ws.Range("T8").Value = 0.02290529633450905
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.406858
$ws.Range("H9").Value = 64.220574
$ws.Range("I9").Value = 0.07340811666748523
$ws.Range("J9").Value = 0.07340811666748523
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.06683666666666
$ws.Range("N9").Value = 87.20050999999998
$ws.Range("O9").Value = 0.1250502234324667
$ws.Range("P9").Value = 0.1250502234324667
$ws.Range("Q9").Value = 622.2296450325265
$ws.Range("R9").Value = 5600.066805292739
$ws.Range("S9").Value = 0.009179701391025615
$ws.Range("T9").Value = 0.009179701391025615
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 261.5073063333334
$ws.Range("H10").Value = 784.521919
$ws.Range("I10").Value = 0.89675742478028
$ws.Range("J10").Value = 0.8967574247802799
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.33695966666667
$ws.Range("N10").Value = 139.010879
$ws.Range("O10").Value = 0.1993490803952133
$ws.Range("P10").Value = 0.1993490803952133
$ws.Range("Q10").Value = 12117.45350610631
$ws.Range("R10").Value = 109057.0815549568
$ws.Range("S10").Value = 0.1787677679675285
$ws.Range("T10").Value = 0.1787677679675284
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 261.5073063333334
$ws.Range("H11").Value = 784.521919
$ws.Range("I11").Value = 0.89675742478028
$ws.Range("J11").Value = 0.8967574247802799
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 84.50960033333332
$ws.Range("N11").Value = 253.528801
$ws.Range("O11").Value = 0.3635739425333109
$ws.Range("P11").Value = 0.3635739425333109
$ws.Range("Q11").Value = 22099.87794247657
$ws.Range("R11").Value = 198898.9014822891
$ws.Range("S11").Value = 0.3260376324233853
$ws.Range("T11").Value = 0.3260376324233853
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 261.5073063333334
$ws.Range("H12").Value = 784.521919
$ws.Range("I12").Value = 0.89675742478028
$ws.Range("J12").Value = 0.8967574247802799
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 72.52790466666666
$ws.Range("N12").Value = 217.583714
$ws.Range("O12").Value = 0.3120267536390091
$ws.Range("P12").Value = 0.3120267536390091
$ws.Range("Q12").Value = 18966.5769833808
$ws.Range("R12").Value = 170699.1928504271
$ws.Range("S12").Value = 0.2798123080558686
$ws.Range("T12").Value = 0.2798123080558686
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 261.5073063333334
$ws.Range("H13").Value = 784.521919
$ws.Range("I13").Value = 0.89675742478028
$ws.Range("J13").Value = 0.8967574247802799
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 29.06683666666666
$ws.Range("N13").Value = 87.20050999999998
$ws.Range("O13").Value = 0.1250502234324667
$ws.Range("P13").Value = 0.1250502234324667
$ws.Range("Q13").Value = 7601.190160330964
$ws.Range("R13").Value = 68410.71144297867
$ws.Range("S13").Value = 0.1121397163334975
$ws.Range("T13").Value = 0.1121397163334975
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 7.296131666666668
$ws.Range("H14").Value = 21.888395
$ws.Range("I14").Value = 0.02501979901057877
$ws.Range("J14").Value = 0.02501979901057877
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 46.33695966666667
$ws.Range("N14").Value = 139.010879
$ws.Range("O14").Value = 0.1993490803952133
$ws.Range("P14").Value = 0.1993490803952133
$ws.Range("Q14").Value = 338.0805587610228
$ws.Range("R14").Value = 3042.725028849205
$ws.Range("S14").Value = 0.004987673924431945
$ws.Range("T14").Value = 0.004987673924431944
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 7.296131666666668
$ws.Range("H15").Value = 21.888395
$ws.Range("I15").Value = 0.02501979901057877
$ws.Range("J15").Value = 0.02501979901057877
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 84.50960033333332
$ws.Range("N15").Value = 253.528801
$ws.Range("O15").Value = 0.3635739425333109
$ws.Range("P15").Value = 0.3635739425333109
$ws.Range("Q15").Value = 616.5931711293772
$ws.Range("R15").Value = 5549.338540164395
$ws.Range("S15").Value = 0.009096546967667152
$ws.Range("T15").Value = 0.009096546967667152
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 7.296131666666668
$ws.Range("H16").Value = 21.888395
$ws.Range("I16").Value = 0.02501979901057877
$ws.Range("J16").Value = 0.02501979901057877
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 72.52790466666666
$ws.Range("N16").Value = 217.583714
$ws.Range("O16").Value = 0.3120267536390091
$ws.Range("P16").Value = 0.3120267536390091
$ws.Range("Q16").Value = 529.1731419554478
$ws.Range("R16").Value = 4762.558277599031
$ws.Range("S16").Value = 0.007806846661971383
$ws.Range("T16").Value = 0.007806846661971383
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 7.296131666666668
$ws.Range("H17").Value = 21.888395
$ws.Range("I17").Value = 0.02501979901057877
$ws.Range("J17").Value = 0.02501979901057877
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 29.06683666666666
$ws.Range("N17").Value = 87.20050999999998
$ws.Range("O17").Value = 0.1250502234324667
$ws.Range("P17").Value = 0.1250502234324667
$ws.Range("Q17").Value = 212.0754674534944
$ws.Range("R17").Value = 1908.67920708145
$ws.Range("S17").Value = 0.003128731456508285
$ws.Range("T17").Value = 0.003128731456508285
